# Add two new "note" worksheets (Revenue note and Income-tax note) at the
# end of the workbook, mirroring the layout of the existing note sheets:
# header row (line_item_0/row_header/value/header_col_0/header_col_1/
# header_col_2/year) followed by data rows of line item / [sub label] /
# value / "$" / year-text / "Consolidated" / year-number.

$wb = $excel.ActiveWorkbook

# Template cells reused so new cells get identical styling without minting
# new style records in styles.xml:
#  - header row formatting (bold + border + centered) lives on row 1 of the
#    first worksheet.
#  - "2023"/"2022" as literal TEXT (not numbers) already exist on that same
#    sheet (E2:E4 = "2023", E5:E7 = "2022"); copying their *values* lets us
#    write the same text without Excel's automatic text->number coercion
#    kicking in (which happens on a plain numeric-looking .Value assignment).
$templateSheet = $wb.Worksheets.Item(1)
$headerTemplate = $templateSheet.Range("A1:G1")
$text2023Template = $templateSheet.Range("E2")
$text2022Template = $templateSheet.Range("E5")

# ---------------------------------------------------------------------
# New sheet 1: "4__c09fe481-4d79-3ba" -> Revenue note
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsRevenue = $wb.Worksheets.Add($null, $lastSheet)
$wsRevenue.Name = "4__c09fe481-4d79-3ba"

$wsRevenue.Range("A1").Value = "line_item_0"
$wsRevenue.Range("B1").Value = "row_header"
$wsRevenue.Range("C1").Value = "value"
$wsRevenue.Range("D1").Value = "header_col_0"
$wsRevenue.Range("E1").Value = "header_col_1"
$wsRevenue.Range("F1").Value = "header_col_2"
$wsRevenue.Range("G1").Value = "year"

$revenueRows = @(
    @("Service revenue", 92279085, 2023),
    @("Minor works",      29612233, 2023),
    @("Revenue",         121891318, 2023),
    @("Service revenue",  85281896, 2022),
    @("Minor works",      24692541, 2022),
    @("Revenue",         109974437, 2022)
)

# New shared strings are interned in first-write order, and the source
# workbook was produced column-by-column (pandas-style export), so column A
# is populated for every data row before moving on to column C, D, etc.
$r = 2
foreach ($row in $revenueRows) {
    $wsRevenue.Range("A$r").Value = $row[0]
    $r = $r + 1
}
$r = 2
foreach ($row in $revenueRows) {
    $wsRevenue.Range("C$r").Value = $row[1]
    $r = $r + 1
}
$r = 2
foreach ($row in $revenueRows) {
    $wsRevenue.Range("D$r").Value = "$"
    $r = $r + 1
}
$r = 2
foreach ($row in $revenueRows) {
    $wsRevenue.Range("F$r").Value = "Consolidated"
    $r = $r + 1
}
$r = 2
foreach ($row in $revenueRows) {
    $wsRevenue.Range("G$r").Value = $row[2]
    $r = $r + 1
}

foreach ($rr in @(2, 3, 4)) {
    $text2023Template.Copy()
    $wsRevenue.Range("E$rr").PasteSpecial(-4163)
}
foreach ($rr in @(5, 6, 7)) {
    $text2022Template.Copy()
    $wsRevenue.Range("E$rr").PasteSpecial(-4163)
}

$headerTemplate.Copy()
$wsRevenue.Range("A1:G1").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# New sheet 2: "6__45029265-bfb6-3df" -> Income tax note
# ---------------------------------------------------------------------
$wsTax = $wb.Worksheets.Add($null, $wsRevenue)
$wsTax.Name = "6__45029265-bfb6-3df"

$wsTax.Range("A1").Value = "line_item_0"
$wsTax.Range("B1").Value = "row_header"
$wsTax.Range("C1").Value = "value"
$wsTax.Range("D1").Value = "header_col_0"
$wsTax.Range("E1").Value = "header_col_1"
$wsTax.Range("F1").Value = "header_col_2"
$wsTax.Range("G1").Value = "year"

$taxRows = @(
    @("Aggregate income tax benefit",                  "Deferred tax - origination and reversal of temporary differences", -481027),
    @("increase in deferred tax assets (note 16)",      "Total", 0),
    @("Loss before income tax benefit",                 "Total", -1684214),
    @("Tax at the statutory tax rate of 30%",           "Total", -505264),
    @("Non-deductible expenses",                        "Total", 24237),
    @("Income tax benefit",                             "Total", -481027)
)

# Same column-major write order as the revenue sheet above, so new shared
# strings are interned in the same sequence as the source export (all of
# column A's new line items, then column B's new sub-label).
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("A$r").Value = $row[0]
    $r = $r + 1
}
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("B$r").Value = $row[1]
    $r = $r + 1
}
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("C$r").Value = $row[2]
    $r = $r + 1
}
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("D$r").Value = "$"
    $r = $r + 1
}
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("F$r").Value = "Consolidated"
    $r = $r + 1
}
$r = 2
foreach ($row in $taxRows) {
    $wsTax.Range("G$r").Value = 2022
    $r = $r + 1
}

foreach ($rr in @(2, 3, 4, 5, 6, 7)) {
    $text2022Template.Copy()
    $wsTax.Range("E$rr").PasteSpecial(-4163)
}

$headerTemplate.Copy()
$wsTax.Range("A1:G1").PasteSpecial(-4122)
